$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new run-log row (row 56) with the same data pattern as the
# other "SKIPPED" rows in the log.
$row = 56
$ws.Cells.Item($row, 1).Value = "2025-08-25 09:41:00 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-25 15:11:00 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($row, 7).Value = 0

# Copy the formatting from the previous row (including the blank, styled
# "Saved PDF" / "Total Rows After" cells) onto the new row so style s="3"
# is applied uniformly across A56:H56, matching the rest of the log.
$ws.Range("A55:H55").Copy()
$ws.Range("A56:H56").PasteSpecial(-4122)
